# Restructure the "Reach" description paragraphs in ReadMe.docx.
$d = $word.ActiveDocument

function Get-ParaByText($txt) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $txt) {
            return $p
        }
    }
    return $null
}

function Set-ParaBodyXml($para, $bodyXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($pkg)
}

$SPACING = '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'

# 1) "The Bureau of Reclamation..." paragraph: add simple spacing pPr, keep run/text.
$p1 = Get-ParaByText "The Bureau of Reclamation is the source/author for the river miles data."
$p1.SpaceAfter = 0
$p1.LineSpacingRule = 0

# 2) Delete the "Reach 6 is:" paragraph entirely (including its paragraph mark).
$p2 = Get-ParaByText "Reach 6 is:"
$p2.Range.Delete()

# 3) "Bernardo HWY 60: 130.62" -> empty paragraph, spacing pPr only, no run.
$p3 = Get-ParaByText "Bernardo HWY 60: 130.62"
$body3 = '<w:body><w:p><w:pPr>' + $SPACING + '</w:pPr></w:p></w:body>'
Set-ParaBodyXml $p3 $body3

# 4) "Bernardo Gage: 130.57" -> "River miles go from high to low going downstream"
$p4 = Get-ParaByText "Bernardo Gage: 130.57"
$body4 = '<w:body><w:p><w:pPr>' + $SPACING + '</w:pPr><w:r><w:t>River miles go from high to low going downstream</w:t></w:r></w:p></w:body>'
Set-ParaBodyXml $p4 $body4

# 5) "SADD: 116.17" -> "Reach 5 is RM 172 – 129 (Isleta to bosque gage)" (3 runs)
$p5 = Get-ParaByText "SADD: 116.17"
$reach5 = '<w:body><w:p><w:pPr>' + $SPACING + '</w:pPr>' + `
  '<w:r><w:t>Reach 5 is RM 172 ' + [char]0x2013 + ' 129 (</w:t></w:r>' + `
  '<w:r><w:t>I</w:t></w:r>' + `
  '<w:r><w:t>sleta to bosque gage)</w:t></w:r>' + `
'</w:p></w:body>'
Set-ParaBodyXml $p5 $reach5

# 6) "SA gage: 115.96" -> Reach 6, Reach 7, Reach 8 paragraphs (3 new paragraphs replacing 1).
$p6 = Get-ParaByText "SA gage: 115.96"
$reach678 = '<w:body>' + `
  '<w:p><w:pPr>' + $SPACING + '</w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Reach 6 is RM </w:t></w:r>' + `
    '<w:r><w:t>130 ' + [char]0x2013 + ' 116 (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Bosuqe</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> to San Acacia gage)</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr>' + $SPACING + '</w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Reach 7 is RM 117 ' + [char]0x2013 + ' 68 (San Acacia to San </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Marical</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> gage)</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr>' + $SPACING + '</w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Reach 8 is RM 67 </w:t></w:r>' + `
    '<w:r><w:t>' + [char]0x2013 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">54 (San </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Marcial</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> gage to Elephant Butte</w:t></w:r>' + `
  '</w:p>' + `
'</w:body>'
Set-ParaBodyXml $p6 $reach678

# 9) Final empty paragraph (previously self-closed <w:p/>) -> add spacing pPr only.
$pLast = $d.Paragraphs.Last
$pLast.SpaceAfter = 0
$pLast.LineSpacingRule = 0

Write-Host "Edit complete"
